# Fruta / hortaliza, semanal
# Insert a new weekly price-report row above row 92 (Locoto / Segunda,
# week of 2022-06-10), pushing the existing rows 92:110 down to 93:111.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 92; Excel shifts rows 92:110 down to 93:111
# and the sheet's used-range dimension grows to A1:R111 automatically.
$ws.Rows.Item(92).Insert()

# Populate the new row 92 with the new week's record.
$ws.Cells.Item(92, 1).Value = 1
$ws.Cells.Item(92, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(92, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(92, 4).Value = 44722
$ws.Cells.Item(92, 5).Value = 15
$ws.Cells.Item(92, 6).Value = 100112042
$ws.Cells.Item(92, 7).Value = "Locoto"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Segunda"
$ws.Cells.Item(92, 10).Value = 150
$ws.Cells.Item(92, 11).Value = 18000
$ws.Cells.Item(92, 12).Value = 20000
$ws.Cells.Item(92, 13).Value = 19000
$ws.Cells.Item(92, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(92, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(92, 16).Value = 950
$ws.Cells.Item(92, 17).Value = 20
$ws.Cells.Item(92, 18).Value = "Hortaliza"
